# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型" sheets
# to reflect the latest scrape output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 41
$ws1.Range("F4").Value = 34
$ws1.Range("F5").Value = 56
$ws1.Range("F6").Value = 55
$ws1.Range("F8").Value = 3793
$ws1.Range("F9").Value = 81
$ws1.Range("F10").Value = 4477
$ws1.Range("F12").Value = 1121
$ws1.Range("F13").Value = 60

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 41
$ws4.Range("F4").Value = 34
$ws4.Range("F5").Value = 56
$ws4.Range("F6").Value = 55
$ws4.Range("F9").Value = 3793
$ws4.Range("F10").Value = 81
$ws4.Range("F11").Value = 4477
$ws4.Range("F13").Value = 1121
$ws4.Range("F14").Value = 60
